# Auto-generated PowerShell Excel COM-interop script
# Applies the Liga Classica dataset update: adds 7 new teams and
# extends sheets 'Geral', 'Turno 2', 'Classif Turno 2', and the six
# 'Mes - ...' sheets from 10 to 17 teams (rows 2-18).

$wb = $excel.ActiveWorkbook

## --- Sheet 'Geral' (sheet1): rows 2-18, columns A:AM ---
$ws1 = $wb.Worksheets.Item("Geral")
$ws1.Range("A2:AM2").Copy() | Out-Null
$ws1.Range("A12:AM12").PasteSpecial(-4122) | Out-Null
$ws1.Range("A13:AM13").PasteSpecial(-4122) | Out-Null
$ws1.Range("A14:AM14").PasteSpecial(-4122) | Out-Null
$ws1.Range("A15:AM15").PasteSpecial(-4122) | Out-Null
$ws1.Range("A16:AM16").PasteSpecial(-4122) | Out-Null
$ws1.Range("A17:AM17").PasteSpecial(-4122) | Out-Null
$ws1.Range("A18:AM18").PasteSpecial(-4122) | Out-Null
$ws1.Application.CutCopyMode = 0
$ws1.Range("A2").Value2 = 'bugredasmissões'
$ws1.Range("B2:AM2").Value2 = 0
$ws1.Range("A3").Value2 = 'C R Juvenal'
$ws1.Range("B3:AM3").Value2 = 0
$ws1.Range("A4").Value2 = 'Doug Leal F.C'
$ws1.Range("B4:AM4").Value2 = 0
$ws1.Range("A5").Value2 = 'Esquadrão Gazembrino'
$ws1.Range("B5:AM5").Value2 = 0
$ws1.Range("A6").Value2 = 'FBC Colorado'
$ws1.Range("B6:AM6").Value2 = 0
$ws1.Range("A7").Value2 = 'GaúchoDaFronteira F.C'
$ws1.Range("B7:AM7").Value2 = 0
$ws1.Range("A8").Value2 = 'GE Bebum'
$ws1.Range("B8:AM8").Value2 = 0
$ws1.Range("A9").Value2 = 'Grêmio_Campeão_LA_27'
$ws1.Range("B9:AM9").Value2 = 0
$ws1.Range("A10").Value2 = 'JV5 Tricolor Gaúcho'
$ws1.Range("B10:AM10").Value2 = 0
$ws1.Range("A11").Value2 = 'La Primeira Patada Es Nuestra'
$ws1.Range("B11:AM11").Value2 = 0
$ws1.Range("A12").Value2 = 'lsauer fc'
$ws1.Range("B12:AM12").Value2 = 0
$ws1.Range("A13").Value2 = 'Medonho´s F.C.'
$ws1.Range("B13:AM13").Value2 = 0
$ws1.Range("A14").Value2 = 'NHU PORÃ SAF.'
$ws1.Range("B14:AM14").Value2 = 0
$ws1.Range("A15").Value2 = 'Pontaç0 F.C.'
$ws1.Range("B15:AM15").Value2 = 0
$ws1.Range("A16").Value2 = 'SC 100 Sono'
$ws1.Range("B16:AM16").Value2 = 0
$ws1.Range("A17").Value2 = 'SC ÉoINTER!'
$ws1.Range("B17:AM17").Value2 = 0
$ws1.Range("A18").Value2 = 'Texas Club 2026'
$ws1.Range("B18:AM18").Value2 = 0

## --- Sheet 'Turno 2' (sheet2): rows 2-18, columns A:T ---
$ws2 = $wb.Worksheets.Item("Turno 2")
$ws2.Range("A2:T2").Copy() | Out-Null
$ws2.Range("A12:T12").PasteSpecial(-4122) | Out-Null
$ws2.Range("A13:T13").PasteSpecial(-4122) | Out-Null
$ws2.Range("A14:T14").PasteSpecial(-4122) | Out-Null
$ws2.Range("A15:T15").PasteSpecial(-4122) | Out-Null
$ws2.Range("A16:T16").PasteSpecial(-4122) | Out-Null
$ws2.Range("A17:T17").PasteSpecial(-4122) | Out-Null
$ws2.Range("A18:T18").PasteSpecial(-4122) | Out-Null
$ws2.Application.CutCopyMode = 0
$ws2.Range("A2").Value2 = 'bugredasmissões'
$ws2.Range("B2:T2").Value2 = 0
$ws2.Range("A3").Value2 = 'C R Juvenal'
$ws2.Range("B3:T3").Value2 = 0
$ws2.Range("A4").Value2 = 'Doug Leal F.C'
$ws2.Range("B4:T4").Value2 = 0
$ws2.Range("A5").Value2 = 'Esquadrão Gazembrino'
$ws2.Range("B5:T5").Value2 = 0
$ws2.Range("A6").Value2 = 'FBC Colorado'
$ws2.Range("B6:T6").Value2 = 0
$ws2.Range("A7").Value2 = 'GaúchoDaFronteira F.C'
$ws2.Range("B7:T7").Value2 = 0
$ws2.Range("A8").Value2 = 'GE Bebum'
$ws2.Range("B8:T8").Value2 = 0
$ws2.Range("A9").Value2 = 'Grêmio_Campeão_LA_27'
$ws2.Range("B9:T9").Value2 = 0
$ws2.Range("A10").Value2 = 'JV5 Tricolor Gaúcho'
$ws2.Range("B10:T10").Value2 = 0
$ws2.Range("A11").Value2 = 'La Primeira Patada Es Nuestra'
$ws2.Range("B11:T11").Value2 = 0
$ws2.Range("A12").Value2 = 'lsauer fc'
$ws2.Range("B12:T12").Value2 = 0
$ws2.Range("A13").Value2 = 'Medonho´s F.C.'
$ws2.Range("B13:T13").Value2 = 0
$ws2.Range("A14").Value2 = 'NHU PORÃ SAF.'
$ws2.Range("B14:T14").Value2 = 0
$ws2.Range("A15").Value2 = 'Pontaç0 F.C.'
$ws2.Range("B15:T15").Value2 = 0
$ws2.Range("A16").Value2 = 'SC 100 Sono'
$ws2.Range("B16:T16").Value2 = 0
$ws2.Range("A17").Value2 = 'SC ÉoINTER!'
$ws2.Range("B17:T17").Value2 = 0
$ws2.Range("A18").Value2 = 'Texas Club 2026'
$ws2.Range("B18:T18").Value2 = 0

## --- Sheets 'Classif Turno 2' + monthly sheets: rows 2-18, columns A:B ---
$monthSheetNames = @(
    'Classif Turno 2',
    'Mês - Janeiro',
    'Mês - Fevereiro',
    'Mês - Março',
    'Mês - Abril',
    'Mês - Maio',
    'Mês - Julho'
)

foreach ($sheetName in $monthSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A2:B2").Copy() | Out-Null
    $ws.Range("A12:B12").PasteSpecial(-4122) | Out-Null
    $ws.Range("A13:B13").PasteSpecial(-4122) | Out-Null
    $ws.Range("A14:B14").PasteSpecial(-4122) | Out-Null
    $ws.Range("A15:B15").PasteSpecial(-4122) | Out-Null
    $ws.Range("A16:B16").PasteSpecial(-4122) | Out-Null
    $ws.Range("A17:B17").PasteSpecial(-4122) | Out-Null
    $ws.Range("A18:B18").PasteSpecial(-4122) | Out-Null
    $ws.Application.CutCopyMode = 0
    $ws.Range("A2").Value2 = 'bugredasmissões'
    $ws.Range("B2").Value2 = 0
    $ws.Range("A3").Value2 = 'La Primeira Patada Es Nuestra'
    $ws.Range("B3").Value2 = 0
    $ws.Range("A4").Value2 = 'SC ÉoINTER!'
    $ws.Range("B4").Value2 = 0
    $ws.Range("A5").Value2 = 'SC 100 Sono'
    $ws.Range("B5").Value2 = 0
    $ws.Range("A6").Value2 = 'Pontaç0 F.C.'
    $ws.Range("B6").Value2 = 0
    $ws.Range("A7").Value2 = 'NHU PORÃ SAF.'
    $ws.Range("B7").Value2 = 0
    $ws.Range("A8").Value2 = 'Medonho´s F.C.'
    $ws.Range("B8").Value2 = 0
    $ws.Range("A9").Value2 = 'lsauer fc'
    $ws.Range("B9").Value2 = 0
    $ws.Range("A10").Value2 = 'JV5 Tricolor Gaúcho'
    $ws.Range("B10").Value2 = 0
    $ws.Range("A11").Value2 = 'C R Juvenal'
    $ws.Range("B11").Value2 = 0
    $ws.Range("A12").Value2 = 'Grêmio_Campeão_LA_27'
    $ws.Range("B12").Value2 = 0
    $ws.Range("A13").Value2 = 'GE Bebum'
    $ws.Range("B13").Value2 = 0
    $ws.Range("A14").Value2 = 'GaúchoDaFronteira F.C'
    $ws.Range("B14").Value2 = 0
    $ws.Range("A15").Value2 = 'FBC Colorado'
    $ws.Range("B15").Value2 = 0
    $ws.Range("A16").Value2 = 'Esquadrão Gazembrino'
    $ws.Range("B16").Value2 = 0
    $ws.Range("A17").Value2 = 'Doug Leal F.C'
    $ws.Range("B17").Value2 = 0
    $ws.Range("A18").Value2 = 'Texas Club 2026'
    $ws.Range("B18").Value2 = 0
}

Write-Host "Liga Classica datasets updated."
